# The underlying data rows for "Top Ten Weather Stories" (row 2) and
# "Canada's top ten weather stories of 2013" (row 3) swap places:
#   - A2/A3 (title text) swap
#   - E2/E3 (uri text) swap
# B/C/D are identical between the two rows already, so nothing to do there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = $ws.Range("A2").Value()
$a3 = $ws.Range("A3").Value()
$e2 = $ws.Range("E2").Value()
$e3 = $ws.Range("E3").Value()

$ws.Range("A2").Value = $a3
$ws.Range("A3").Value = $a2
$ws.Range("E2").Value = $e3
$ws.Range("E3").Value = $e2
